$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Insert a new row at row 3 (shifts the existing rows 3:10 down to 4:11,
# carrying their values/formulas/styles with them)
$ws.Rows.Item(3).Insert()

# Fill in the new row with the CRM opened 20210221 data point
$ws.Cells.Item(3, 1).Value = 20210221
$ws.Cells.Item(3, 2).Value = 2215.554
$ws.Cells.Item(3, 3).Value = 2234.0700000000002
$ws.Cells.Item(3, 4).Formula = "=100*(B3-C3)/C3"
$ws.Cells.Item(3, 5).Value = 141
$ws.Cells.Item(3, 6).Value = "CRM opened 20210221"

# Match the workbook's last-used selection
$ws.Range("F4").Select()
